$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.942.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.03%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.791.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.16%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'358.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.85%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'109.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.38%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.85%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +0.06%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -1.93%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'40.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.95%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +2.08%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.75%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'19.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.38%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'7.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.94%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.228.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.54%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "'WrappedEther"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'2.796.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.50%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "'Polygon"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.947"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.88%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'51.887.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.02%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'7.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.06%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -2.05%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'13.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.35%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  -1.61%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'270.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.72%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'70.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.23%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -1.60%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'26.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.09%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -0.04%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +18.15%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'10.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.77%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'2.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.31%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'52.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.76%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'Filecoin"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'6.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.82%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "'InjectiveProtocol"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'34.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.97%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "'VeChain"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'0.0464"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.83%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.0850"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.57%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -4.04%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +0.03%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'18.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.30%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -2.35%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -3.58%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'2.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.17%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.114"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.90%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -2.19%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'119.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.11%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'21.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -6.09%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.080.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.85%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -2.22%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'2.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.17%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'5.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.39%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.954"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.04%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +30.59%  "
$ws.Range("E51").Style = "Normal"

